$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fill in the previously blank "Actual Result" / expected-result cells for
# the "moving down" (top-of-canvas) row with the newly-added text.
$text = "player stops moving down and is able to move along top of canvas"
$ws.Range("C8").Value = $text
$ws.Range("D8").Value = $text

# Move the selection to C11 (matches the saved view state in the diff).
$ws.Range("C11").Select()
